# Atualizado por script em 05-11-2023 08:45
#
# The match-result/odds columns (F:V) for several rows got reshuffled
# (rows were re-ordered upstream while the row index (A) and match date
# (E) stayed put), and six brand-new fixtures were appended at the end
# of the sheet (rows 82-87).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Rows {
    param($ws, [hashtable]$mapping)

    # Snapshot every source row's F:V payload BEFORE any writes happen,
    # so overlapping cycles don't clobber data we still need to read.
    $snapshots = @{}
    foreach ($srcRow in $mapping.Values) {
        if (-not $snapshots.ContainsKey($srcRow)) {
            $snapshots[$srcRow] = $ws.Range("F$srcRow`:V$srcRow").Value2
        }
    }

    foreach ($newRow in $mapping.Keys) {
        $srcRow = $mapping[$newRow]
        $ws.Range("F$newRow`:V$newRow").Value2 = $snapshots[$srcRow]
    }
}

# Block 1: rows 15-17 (3-cycle)
Rotate-Rows $ws @{ 15 = 16; 16 = 17; 17 = 15 }

# Block 2: rows 19-23 (5-cycle)
Rotate-Rows $ws @{ 19 = 22; 20 = 19; 21 = 23; 22 = 20; 23 = 21 }

# Block 3: rows 38-41 (4-cycle)
Rotate-Rows $ws @{ 38 = 41; 39 = 38; 40 = 39; 41 = 40 }

# Block 4: rows 54-56 (3-cycle)
Rotate-Rows $ws @{ 54 = 56; 55 = 54; 56 = 55 }

# Block 5: rows 68-69 (swap)
Rotate-Rows $ws @{ 68 = 69; 69 = 68 }

# Block 6: rows 75-76 (swap)
Rotate-Rows $ws @{ 75 = 76; 76 = 75 }

# Six brand-new fixtures appended as rows 82-87 (Indice 81-86).
# Carry the existing formatting (bold/border/centre for the index column,
# date-time number format for the match-date column) down onto the new
# rows, same as every prior row in the sheet.
$ws.Range("A81").Copy()
$ws.Range("A82:A87").PasteSpecial(-4122)
$ws.Range("E81").Copy()
$ws.Range("E82:E87").PasteSpecial(-4122)

$newRows = @(
    @{ Row=82; A=81; E=45234.64583333334; F="Cardiff Metropolitan"; G=3; H="Connahs Q.";   I=1; J=4.39; K="02/11/2023 08:13"; L=5.35; M="04/11/2023 15:23"; N=3.87; O="02/11/2023 08:13"; P=4.1;  Q="04/11/2023 15:23"; R=1.62; S="02/11/2023 08:13"; T=1.6;  U="04/11/2023 15:23"; V="https://www.betexplorer.com/football/wales/cymru-premier/cardiff-metropolitan-university-connahs-q/hl0qMU0q/" },
    @{ Row=83; A=82; E=45234.64583333334; F="Caernarfon";           G=2; H="Penybont";      I=4; J=2.62; K="02/11/2023 08:13"; L=2.62; M="04/11/2023 15:22"; N=3.34; O="02/11/2023 08:13"; P=3.77; Q="04/11/2023 15:22"; R=2.38; S="02/11/2023 08:13"; T=2.44; U="04/11/2023 15:22"; V="https://www.betexplorer.com/football/wales/cymru-premier/caernarfon-penybont/0d8HQAhS/" },
    @{ Row=84; A=83; E=45234.64583333334; F="Haverfordwest";       G=5; H="Colwyn Bay";     I=0; J=1.69; K="02/11/2023 08:13"; L=1.83; M="04/11/2023 15:22"; N=3.75; O="02/11/2023 08:13"; P=3.64; Q="04/11/2023 15:24"; R=4.07; S="02/11/2023 08:13"; T=4.26; U="04/11/2023 15:24"; V="https://www.betexplorer.com/football/wales/cymru-premier/haverfordwest-colwyn-bay/WvamLlGk/" },
    @{ Row=85; A=84; E=45234.64583333334; F="Bala";                G=0; H="Pontypridd";    I=0; J=1.6;  K="02/11/2023 08:13"; L=1.57; M="04/11/2023 15:22"; N=3.75; O="02/11/2023 08:13"; P=3.68; Q="04/11/2023 15:22"; R=4.72; S="02/11/2023 08:13"; T=6.74; U="04/11/2023 15:22"; V="https://www.betexplorer.com/football/wales/cymru-premier/bala-pontypridd-united/vsK8SWNF/" },
    @{ Row=86; A=85; E=45234.64583333334; F="Barry";               G=0; H="Aberystwyth";   I=1; J=1.74; K="02/11/2023 08:13"; L=1.69; M="04/11/2023 14:58"; N=3.68; O="02/11/2023 08:13"; P=4.07; Q="04/11/2023 14:58"; R=3.91; S="02/11/2023 08:13"; T=4.5;  U="04/11/2023 14:58"; V="https://www.betexplorer.com/football/wales/cymru-premier/barry-town-aberystwyth/M59DRjwM/" },
    @{ Row=87; A=86; E=45234.64583333334; F="Newtown";             G=0; H="TNS";           I=2; J=7.29; K="02/11/2023 08:13"; L=8.68; M="04/11/2023 15:13"; N=5.68; O="02/11/2023 08:13"; P=5.96; Q="04/11/2023 15:13"; R=1.27; S="02/11/2023 08:13"; T=1.29; U="04/11/2023 15:10"; V="https://www.betexplorer.com/football/wales/cymru-premier/newtown-tns/E3lhK8Ve/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.A          # Indice
    $ws.Cells.Item($r, 2).Value = "wales"        # pais
    $ws.Cells.Item($r, 3).Value = "cymru-premier"# torneio
    $ws.Cells.Item($r, 4).Value = "2023-2024"    # temporada
    $ws.Cells.Item($r, 5).Value = $nr.E          # data_partida
    $ws.Cells.Item($r, 6).Value = $nr.F          # home
    $ws.Cells.Item($r, 7).Value = $nr.G          # home_ft_gols
    $ws.Cells.Item($r, 8).Value = $nr.H          # away
    $ws.Cells.Item($r, 9).Value = $nr.I          # away_ft_gols
    $ws.Cells.Item($r, 10).Value = $nr.J         # home_opening_odds
    $ws.Cells.Item($r, 11).Value = $nr.K         # home_opening_data_hora
    $ws.Cells.Item($r, 12).Value = $nr.L         # home_closing_odds
    $ws.Cells.Item($r, 13).Value = $nr.M         # home_closing_data_hora
    $ws.Cells.Item($r, 14).Value = $nr.N         # draw_opening_odds
    $ws.Cells.Item($r, 15).Value = $nr.O         # draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value = $nr.P         # draw_closing_odds
    $ws.Cells.Item($r, 17).Value = $nr.Q         # draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value = $nr.R         # away_opening_odds
    $ws.Cells.Item($r, 19).Value = $nr.S         # away_opening_data_hora
    $ws.Cells.Item($r, 20).Value = $nr.T         # away_closing_odds
    $ws.Cells.Item($r, 21).Value = $nr.U         # away_closing_data_hora
    $ws.Cells.Item($r, 22).Value = $nr.V         # url_partida
}
